# Mark Cleolia Jasmine Nasution (row 8) as having completed Tugas 2, 3, and 4,
# and mark Nouridza Juniansah Ridhan (row 27) as having completed Tugas 5.
# Completed tasks are shown with a Wingdings check-mark glyph ("ü"), matching
# the formatting already used by every other completed cell in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$checkCells = @("D8", "E8", "F8", "G27")

foreach ($addr in $checkCells) {
    $cell = $ws.Range($addr)
    $cell.Value = "ü"
    $cell.Font.Name = "Wingdings"
}
